$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Schedule")
$ws2 = $wb.Worksheets.Item("Detailed")

# ------------------------------------------------------------------
# Sheet "Schedule": rows 2-3 updated, new row 4 inserted with values
# ------------------------------------------------------------------

# Row 2 updates
$ws1.Range("B2").Value = 46039.8125
$ws1.Range("C2").Value = 12.5
$ws1.Range("D2").Value = 47.25
$ws1.Range("E2").Value = 652.71611925
$ws1.Range("F2").Value = 13.81409776190476

# Row 3 updates (old row 3 data shifts into the now-earlier slot)
$ws1.Range("A3").Value = 46039.85416666666
$ws1.Range("B3").Value = 46040.02083333334
$ws1.Range("C3").Value = 4
$ws1.Range("D3").Value = 15.12
$ws1.Range("E3").Value = 296.949861
$ws1.Range("F3").Value = 19.63954107142857

# New row 4 (previously-last schedule entry, now split into rows 3 & 4)
$ws1.Range("A4").Value = 46040.3125
$ws1.Range("B4").Value = 46040.79166666666
$ws1.Range("C4").Value = 11.5
$ws1.Range("D4").Value = 43.47
$ws1.Range("E4").Value = 71.07485774999998
$ws1.Range("F4").Value = 1.635032384403036

# Apply the same date/time number format used by the other Start/Stop Time
# cells (style index 2) to the newly added row 4 cells
$ws1.Range("A4:B4").NumberFormat = $ws1.Range("A2").NumberFormat

# ------------------------------------------------------------------
# Sheet "Detailed": price (B) and status (C/E) updates for rows 36-94
# ------------------------------------------------------------------

$ws2.Range("B36").Value = -0.57248
$ws2.Range("B37").Value = -2.97056
$ws2.Range("B38").Value = -2.80872
$ws2.Range("B39").Value = -3.69268
$ws2.Range("C39").Value = "historical"
$ws2.Range("B40").Value = 17.74626
$ws2.Range("C40").Value = "historical"
$ws2.Range("B41").Value = 57.18142
$ws2.Range("E41").Value = "OFF"
$ws2.Range("B42").Value = 56.98
$ws2.Range("E42").Value = "OFF"
$ws2.Range("B43").Value = 46.39479
$ws2.Range("B44").Value = 30.53091
$ws2.Range("B45").Value = 36.2
$ws2.Range("B46").Value = 36.05933
$ws2.Range("B47").Value = 36.2
$ws2.Range("B48").Value = 47.05815
$ws2.Range("E48").Value = "ON"
$ws2.Range("B49").Value = 36.0604
$ws2.Range("E49").Value = "ON"
$ws2.Range("E50").Value = "ON"
$ws2.Range("B52").Value = 31.02048
$ws2.Range("B53").Value = 36.2
$ws2.Range("B56").Value = 36.2
$ws2.Range("B57").Value = 56.97996
$ws2.Range("B58").Value = 50.31057
$ws2.Range("B60").Value = 56.98
$ws2.Range("B61").Value = 56.98
$ws2.Range("B62").Value = 56.98
$ws2.Range("B64").Value = 36.0595
$ws2.Range("E64").Value = "OFF"
$ws2.Range("B65").Value = 28.67846
$ws2.Range("B66").Value = 0.62605
$ws2.Range("B67").Value = 19.29447
$ws2.Range("B70").Value = 26.41845
$ws2.Range("B71").Value = 22.07
$ws2.Range("B72").Value = 0.7
$ws2.Range("B73").Value = 20.59504
$ws2.Range("B74").Value = 0.51
$ws2.Range("B75").Value = 0.64597
$ws2.Range("B76").Value = 0.7
$ws2.Range("B77").Value = 0.00025
$ws2.Range("B78").Value = -4.9548
$ws2.Range("B79").Value = -6.70389
$ws2.Range("B80").Value = -5.41349
$ws2.Range("B81").Value = -6.85668
$ws2.Range("B82").Value = -6.73561
$ws2.Range("B83").Value = -6.89381
$ws2.Range("B84").Value = -7.08241
$ws2.Range("B85").Value = -1.8149
$ws2.Range("B86").Value = -1.41518
$ws2.Range("B87").Value = 0.00975
$ws2.Range("B89").Value = 46.5704
$ws2.Range("B90").Value = 56.98
$ws2.Range("B91").Value = 55.47778
$ws2.Range("B92").Value = 46.8495
$ws2.Range("B94").Value = 56.0398
